# CaseFiling.xlsx - Hawaii Case Filing Juvenile information update
#
# Summary of the intended edit (reverse-engineered from the OOXML diff):
#   1. Sheet renamed from "Prosecution Case Filing" to "Case Filing Decision"
#   2. The title cell (A1) text "Hawaii PA Document" is replaced with
#      "Hawaii PA Document (JJIS)" (this is the only real shared-string
#      content change; the other "Hawaii PA Document" occurrence used
#      elsewhere in the workbook is untouched because D106 kept its
#      original "Employee Occupation Text" value, it only shifted index)
#   3. Three label cells (A106 "Occupation", A114 "Relation to Juvenile",
#      A115 "Employer") are highlighted in bold red to flag them for
#      review/attention
#   4. The window/selection was scrolled so rows 73-74 are visible with
#      A73:A74 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the worksheet tab
$ws.Name = "Case Filing Decision"

# 2) Update the report title cell - triggers the shared-string table to
#    drop the now-unused "Hawaii PA Document" string and append the new
#    "Hawaii PA Document (JJIS)" string, matching the target workbook.
$ws.Range("A1").Value = "Hawaii PA Document (JJIS)"

# 3) Flag the three relevant rows with a bold red font
$ws.Range("A106").Font.Bold = $true
$ws.Range("A106").Font.Color = 255

$ws.Range("A114").Font.Bold = $true
$ws.Range("A114").Font.Color = 255

$ws.Range("A115").Font.Bold = $true
$ws.Range("A115").Font.Color = 255

# 4) Update the visible selection/scroll position to rows 73-74
$excel.Goto($ws.Range("A73:A74"), $true)

# 5) Best-effort: nudge the workbook window position to match the saved
#    view (xWindow/yWindow). Not all hosts persist this back to the
#    workbookView element, but we still set it for parity when supported.
$win = $excel.ActiveWindow
$win.Left = 5500
$win.Top = 1720
